$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Price" cells whose new values look like plain numbers to stay
# stored as text (matching the original inline-string cell type) by applying
# a text number format before writing the value.
$textPriceCells = @(
    'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D13', 'D14', 'D16', 'D17', 'D19', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D45', 'D47', 'D48', 'D49', 'D50', 'D51'
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated coin name / link / price / volume values.
$ws.Range('D2').Value = '27.559.88'
$ws.Range('E2').Value = '  +2.26%  '
$ws.Range('D3').Value = '1.790.26'
$ws.Range('E3').Value = '  +4.06%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '313.61'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = '0.5376'
$ws.Range('E7').Value = '  +10.85%  '
$ws.Range('D8').Value = '0.3764'
$ws.Range('E8').Value = '  +7.67%  '
$ws.Range('D9').Value = '42.86'
$ws.Range('E9').Value = '  +1.81%  '
$ws.Range('D10').Value = '0.07510'
$ws.Range('E10').Value = '  +3.44%  '
$ws.Range('E11').Value = '  +6.27%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = '20.88'
$ws.Range('E13').Value = '  +4.96%  '
$ws.Range('D14').Value = '6.167'
$ws.Range('E14').Value = '  +5.24%  '
$ws.Range('D15').Value = '1.788.57'
$ws.Range('E15').Value = '  +3.44%  '
$ws.Range('D16').Value = '7.072'
$ws.Range('E16').Value = '  +3.16%  '
$ws.Range('D17').Value = '90.82'
$ws.Range('E17').Value = '  +4.76%  '
$ws.Range('D19').Value = '0.06495'
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '16.97'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').Value = '5.938'
$ws.Range('D23').Value = '27.595.00'
$ws.Range('E23').Value = '  +2.21%  '
$ws.Range('D24').Value = '11.21'
$ws.Range('E24').Value = '  +3.71%  '
$ws.Range('D25').Value = '2.085'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '20.48'
$ws.Range('E26').Value = '  +2.76%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '155.38'
$ws.Range('E27').Value = '  +0.79%  '
$ws.Range('D28').Value = '2.381'
$ws.Range('E28').Value = '  +15.60%  '
$ws.Range('D29').Value = '1.995.84'
$ws.Range('E29').Value = '  +3.75%  '
$ws.Range('D30').Value = '121.83'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').Value = '1.122'
$ws.Range('E31').Value = '  +9.23%  '
$ws.Range('D32').Value = '0.1031'
$ws.Range('E32').Value = '  +10.84%  '
$ws.Range('D33').Value = '5.669'
$ws.Range('E33').Value = '  +5.72%  '
$ws.Range('D34').Value = '3.599'
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('D35').Value = '0.02287'
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('D36').Value = '8.724'
$ws.Range('E36').Value = '  +15.99%  '
$ws.Range('D37').Value = '0.06023'
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2089'
$ws.Range('E38').Value = '  +4.70%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = '4.989'
$ws.Range('E39').Value = '  +5.07%  '
$ws.Range('D40').Value = '11.40'
$ws.Range('E40').Value = '  +3.82%  '
$ws.Range('D41').Value = '0.6241'
$ws.Range('E41').Value = '  +4.36%  '
$ws.Range('D42').Value = '1.412'
$ws.Range('E42').Value = '  -1.88%  '
$ws.Range('D43').Value = '0.9997'
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +5.20%  '
$ws.Range('D45').Value = '13.34'
$ws.Range('E45').Value = '  +4.66%  '
$ws.Range('E46').Value = '  +4.33%  '
$ws.Range('D47').Value = '3.634'
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').Value = '121.48'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('D49').Value = '1.911'
$ws.Range('E49').Value = '  +4.08%  '
$ws.Range('D50').Value = '1.133'
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('D51').Value = '0.06747'
$ws.Range('E51').Value = '  +1.51%  '
